# Updates the cryptos table (columns B-E) to the latest scraped values.
# A leading apostrophe forces Excel to store numeric-looking "Price" strings
# (e.g. "212.30") as literal text instead of auto-converting them to numbers,
# matching the source data which always stores these columns as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    $ws.Range($addr).Value = $text
}

# Row 2
Set-TextCell "D2" '86.929.13'
Set-TextCell "E2" '  +5.56%  '

# Row 3
Set-TextCell "D3" '3.287.51'
Set-TextCell "E3" '  +3.55%  '

# Row 4
Set-TextCell "E4" '  -0.10%  '

# Row 5
Set-TextCell "D5" '''212.30'
Set-TextCell "E5" '  -2.50%  '

# Row 6
Set-TextCell "D6" '''629.25'
Set-TextCell "E6" '  +1.52%  '

# Row 7
Set-TextCell "D7" '''0.379'
Set-TextCell "E7" '  +30.65%  '

# Row 8
Set-TextCell "D8" '''0.666'
Set-TextCell "E8" '  +14.50%  '

# Row 9
Set-TextCell "E9" '  -0.03%  '

# Row 10
Set-TextCell "D10" '3.286.74'
Set-TextCell "E10" '  +3.61%  '

# Row 11
Set-TextCell "D11" '''0.580'
Set-TextCell "E11" '  -2.42%  '

# Row 12
Set-TextCell "E12" '  +7.36%  '

# Row 13
Set-TextCell "D13" '''0.0000259'
Set-TextCell "E13" '  +1.19%  '

# Row 14
Set-TextCell "B14" 'Avalanche'
Set-TextCell "C14" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell "D14" '''34.45'
Set-TextCell "E14" '  +7.29%  '

# Row 15
Set-TextCell "B15" 'WrappedliquidstakedEther2.0'
Set-TextCell "C15" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell "D15" '3.877.06'
Set-TextCell "E15" '  +3.38%  '

# Row 16
Set-TextCell "D16" '''5.30'
Set-TextCell "E16" '  -0.15%  '

# Row 17
Set-TextCell "D17" '86.571.17'
Set-TextCell "E17" '  +5.56%  '

# Row 18
Set-TextCell "D18" '3.281.32'
Set-TextCell "E18" '  +3.38%  '

# Row 19
Set-TextCell "D19" '''14.16'
Set-TextCell "E19" '  +1.43%  '

# Row 20
Set-TextCell "D20" '''3.02'
Set-TextCell "E20" '  -6.78%  '

# Row 21
Set-TextCell "D21" '''431.69'
Set-TextCell "E21" '  -1.42%  '

# Row 22
Set-TextCell "D22" '''8.96'
Set-TextCell "E22" '  +0.52%  '

# Row 23
Set-TextCell "E23" '  +4.07%  '

# Row 24
Set-TextCell "E24" '  -1.57%  '

# Row 25
Set-TextCell "D25" '''12.54'
Set-TextCell "E25" '  +4.86%  '

# Row 26
Set-TextCell "D26" '''5.15'
Set-TextCell "E26" '  -1.47%  '

# Row 27
Set-TextCell "D27" '3.439.53'
Set-TextCell "E27" '  +3.26%  '

# Row 28
Set-TextCell "D28" '''76.12'
Set-TextCell "E28" '  -1.22%  '

# Row 29
Set-TextCell "D29" '''0.0000130'
Set-TextCell "E29" '  +8.33%  '

# Row 30
Set-TextCell "D30" '''0.999'
Set-TextCell "E30" '  -0.02%  '

# Row 31
Set-TextCell "D31" '''0.179'
Set-TextCell "E31" '  +21.50%  '

# Row 32
Set-TextCell "D32" '''0.996'
Set-TextCell "E32" '  -0.39%  '

# Row 33
Set-TextCell "D33" '''8.84'
Set-TextCell "E33" '  -1.99%  '

# Row 34
Set-TextCell "D34" '''549.76'
Set-TextCell "E34" '  -3.57%  '

# Row 35
Set-TextCell "D35" '''1.44'
Set-TextCell "E35" '  -3.65%  '

# Row 36
Set-TextCell "D36" '''1.95'
Set-TextCell "E36" '  -1.45%  '

# Row 37
Set-TextCell "D37" '''7.02'
Set-TextCell "E37" '  +12.29%  '

# Row 38
Set-TextCell "D38" '''0.138'
Set-TextCell "E38" '  -8.62%  '

# Row 39
Set-TextCell "D39" '''22.57'
Set-TextCell "E39" '  -0.33%  '

# Row 40
Set-TextCell "E40" '  -0.06%  '

# Row 41
Set-TextCell "D41" '''21.59'
Set-TextCell "E41" '  +3.69%  '

# Row 42
Set-TextCell "D42" '''0.395'
Set-TextCell "E42" '  -2.18%  '

# Row 43
Set-TextCell "D43" '''2.01'
Set-TextCell "E43" '  -0.38%  '

# Row 44
Set-TextCell "D44" '''2.96'
Set-TextCell "E44" '  -0.13%  '

# Row 45
Set-TextCell "B45" 'USDe'
Set-TextCell "C45" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell "D45" '''1.00'
Set-TextCell "E45" '  +0.07%  '

# Row 46
Set-TextCell "B46" 'Monero'
Set-TextCell "C46" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell "D46" '''157.65'
Set-TextCell "E46" '  -1.49%  '

# Row 47
Set-TextCell "D47" '''179.53'
Set-TextCell "E47" '  -3.38%  '

# Row 48
Set-TextCell "D48" '''44.39'
Set-TextCell "E48" '  -0.48%  '

# Row 49
Set-TextCell "D49" '''1.31'
Set-TextCell "E49" '  -0.72%  '

# Row 50
Set-TextCell "D50" '''4.27'
Set-TextCell "E50" '  +2.06%  '

# Row 51
Set-TextCell "D51" '''0.625'
Set-TextCell "E51" '  -0.48%  '
